$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 754
$ws.Range("I6").Value = 754
$ws.Range("K6").Value = 2262
$ws.Range("M6").Value = -2150
$ws.Range("H19").Value = 4552.5386
$ws.Range("I19").Value = 8588.083000000001
$ws.Range("J19").Value = 1093.5
$ws.Range("K19").Value = 8588.083000000001
$ws.Range("L19").Value = 1093.5
$ws.Range("M19").Value = -8413.083000000001
$ws.Range("N19").Value = -1443.5
$ws.Range("H31").Value = 2634.3333
$ws.Range("I31").Value = 2634.3333
$ws.Range("K31").Value = 7902.999899999999
$ws.Range("M31").Value = -7672.999899999999
$ws.Range("H38").Value = 269.8
$ws.Range("I38").Value = 110.5
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 331.5
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = 40.5
$ws.Range("N38").Value = -8244
$ws.Range("H39").Value = 367.33334
$ws.Range("I39").Value = 354.75
$ws.Range("J39").Value = 392.5
$ws.Range("K39").Value = 1064.25
$ws.Range("L39").Value = 1177.5
$ws.Range("M39").Value = -768.25
$ws.Range("N39").Value = -1769.5
$ws.Range("H113").Value = 3290.853
$ws.Range("I113").Value = 3251.4285
$ws.Range("J113").Value = 3301.074
$ws.Range("K113").Value = 3251.4285
$ws.Range("L113").Value = 3301.074
$ws.Range("M113").Value = 2.571500000000015
$ws.Range("N113").Value = -9809.074000000001
$ws.Range("H116").Value = 154225.86
$ws.Range("I116").Value = 179100.17
$ws.Range("K116").Value = 179100.17
$ws.Range("M116").Value = -175658.17
$ws.Range("H127").Value = 1004.3889
$ws.Range("I127").Value = 613.5
$ws.Range("J127").Value = 1199.8334
$ws.Range("K127").Value = 1840.5
$ws.Range("L127").Value = 3599.5002
$ws.Range("M127").Value = 3119.5
$ws.Range("N127").Value = -13519.5002
$ws.Range("H132").Value = 2143.7231
$ws.Range("I132").Value = 1249.8363
$ws.Range("J132").Value = 7060.1
$ws.Range("K132").Value = 3749.5089
$ws.Range("L132").Value = 21180.3
$ws.Range("M132").Value = -1219.5089
$ws.Range("N132").Value = -26240.3
$ws.Range("H138").Value = 2056.8132
$ws.Range("I138").Value = 938.8889
$ws.Range("J138").Value = 3088.7437
$ws.Range("K138").Value = 2816.6667
$ws.Range("L138").Value = 9266.231100000001
$ws.Range("M138").Value = 2323.3333
$ws.Range("N138").Value = -19546.2311

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2469.4167
$ws.Range("I122").Value = 2177.5757
$ws.Range("J122").Value = 3111.4666
$ws.Range("K122").Value = 6532.7271
$ws.Range("L122").Value = 9334.399800000001
$ws.Range("M122").Value = -4082.7271
$ws.Range("N122").Value = -14234.3998
$ws.Range("H125").Value = 38904.668
$ws.Range("J125").Value = 38904.668
$ws.Range("L125").Value = 38904.668
$ws.Range("N125").Value = -48744.668
$ws.Range("H135").Value = 35172.375
$ws.Range("J135").Value = 35172.375
$ws.Range("L135").Value = 35172.375
$ws.Range("N135").Value = -45312.375

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3299082.5
$ws.Range("I99").Value = 1482167.9
$ws.Range("J99").Value = 6251569
$ws.Range("K99").Value = 1482167.9
$ws.Range("L99").Value = 6251569
$ws.Range("M99").Value = -1480669.9
$ws.Range("N99").Value = -6254565

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2661.5833
$ws.Range("I31").Value = 1971.925
$ws.Range("J31").Value = 4040.9
$ws.Range("K31").Value = 1971.925
$ws.Range("L31").Value = 4040.9
$ws.Range("M31").Value = -1676.925
$ws.Range("N31").Value = -4630.9
$ws.Range("H34").Value = 2661.5833
$ws.Range("I34").Value = 1971.925
$ws.Range("J34").Value = 4040.9
$ws.Range("K34").Value = 1971.925
$ws.Range("L34").Value = 4040.9
$ws.Range("M34").Value = -1769.925
$ws.Range("N34").Value = -4444.9
$ws.Range("H58").Value = 1500.2572
$ws.Range("I58").Value = 800.4167
$ws.Range("J58").Value = 3027.182
$ws.Range("K58").Value = 800.4167
$ws.Range("L58").Value = 3027.182
$ws.Range("M58").Value = -597.4167
$ws.Range("N58").Value = -3433.182
$ws.Range("H60").Value = 11103
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 11103
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11103
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = -12125
$ws.Range("H105").Value = 1167.6923
$ws.Range("I105").Value = 728
$ws.Range("J105").Value = 1442.5
$ws.Range("K105").Value = 728
$ws.Range("L105").Value = 1442.5
$ws.Range("M105").Value = 1019
$ws.Range("N105").Value = -4936.5
$ws.Range("H122").Value = 1213.2632
$ws.Range("I122").Value = 1004.875
$ws.Range("J122").Value = 1364.8182
$ws.Range("K122").Value = 3014.625
$ws.Range("L122").Value = 4094.4546
$ws.Range("M122").Value = -564.625
$ws.Range("N122").Value = -8994.454600000001
$ws.Range("H132").Value = 2214.85
$ws.Range("I132").Value = 1188.25
$ws.Range("J132").Value = 3754.75
$ws.Range("K132").Value = 3564.75
$ws.Range("L132").Value = 11264.25
$ws.Range("M132").Value = -1034.75
$ws.Range("N132").Value = -16324.25
$ws.Range("H134").Value = 1503.8718
$ws.Range("I134").Value = 959.5925999999999
$ws.Range("J134").Value = 2728.5
$ws.Range("K134").Value = 2878.7778
$ws.Range("L134").Value = 8185.5
$ws.Range("M134").Value = -343.7777999999998
$ws.Range("N134").Value = -13255.5
$ws.Range("H136").Value = 1500.2572
$ws.Range("I136").Value = 800.4167
$ws.Range("J136").Value = 3027.182
$ws.Range("K136").Value = 2401.2501
$ws.Range("L136").Value = 9081.545999999998
$ws.Range("M136").Value = 148.7498999999998
$ws.Range("N136").Value = -14181.546

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 933.3333
$ws.Range("I13").Value = 800
$ws.Range("K13").Value = 2400
$ws.Range("M13").Value = -2232

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1552.68
$ws.Range("I122").Value = 1516.0625
$ws.Range("J122").Value = 1617.7778
$ws.Range("K122").Value = 4548.1875
$ws.Range("L122").Value = 4853.3334
$ws.Range("M122").Value = -2098.1875
$ws.Range("N122").Value = -9753.3334
$ws.Range("H126").Value = 2849.7441
$ws.Range("I126").Value = 2370.36
$ws.Range("K126").Value = 7111.08
$ws.Range("M126").Value = -4641.08

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9162.138000000001
$ws.Range("I132").Value = 2390.2
$ws.Range("J132").Value = 24210.889
$ws.Range("K132").Value = 7170.599999999999
$ws.Range("L132").Value = 72632.667
$ws.Range("M132").Value = -4640.599999999999
$ws.Range("N132").Value = -77692.667
$ws.Range("H133").Value = 28163
$ws.Range("J133").Value = 28163
$ws.Range("L133").Value = 28163
$ws.Range("N133").Value = -33223

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 19178.334
$ws.Range("I58").Value = 5585
$ws.Range("J58").Value = 25975
$ws.Range("K58").Value = 5585
$ws.Range("L58").Value = 25975
$ws.Range("M58").Value = -5277
$ws.Range("N58").Value = -26591
$ws.Range("H107").Value = 553.76
$ws.Range("I107").Value = 503.7857
$ws.Range("J107").Value = 617.36365
$ws.Range("K107").Value = 1511.3571
$ws.Range("L107").Value = 1852.09095
$ws.Range("M107").Value = 408.6428999999998
$ws.Range("N107").Value = -5692.09095
$ws.Range("H132").Value = 3072.3044
$ws.Range("I132").Value = 2718.5454
$ws.Range("J132").Value = 3396.5833
$ws.Range("K132").Value = 8155.6362
$ws.Range("L132").Value = 10189.7499
$ws.Range("M132").Value = -5625.6362
$ws.Range("N132").Value = -15249.7499
